$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refreshed data table (6 data rows x 20 columns, A:T) reflecting the
# new TPM-derived values. Columns A-D are text labels (sending cluster,
# ligand symbol, receptor symbol, target cluster); columns E-T are numeric.
$data = @(
    @("ECs",   "Nppa", "Npr2", "ECs",   2, 0.6666666666666666, 0.1893306666666667, 0.5679920000000001, 0.6308162521878971, 0.6308162521878971, 3, 1, 13.65974866666667, 40.979246,         0.383515069756994,  0.3835150697569941, 2.586209321559112, 23.275883894032, 0.2419275389616869, 0.2419275389616869),
    @("ECs",   "Nppa", "Npr2", "FAPs",  2, 0.6666666666666666, 0.1893306666666667, 0.5679920000000001, 0.6308162521878971, 0.6308162521878971, 3, 1, 17.104695,         51.314085,         0.4802363832729211, 0.4802363832729211, 3.23844330748,     29.14598976732,  0.3029409154604946, 0.3029409154604946),
    @("ECs",   "Nppa", "Npr2", "MuSCs", 2, 0.6666666666666666, 0.1893306666666667, 0.5679920000000001, 0.6308162521878971, 0.6308162521878971, 3, 1, 4.852797333333333, 14.558392,         0.1362485469700849, 0.1362485469700849, 0.9187833543182222, 8.269050188864, 0.08594779776571561, 0.08594779776571561),
    @("MuSCs", "Nppa", "Npr2", "ECs",   2, 0.6666666666666666, 0.1108053333333333, 0.332416,           0.3691837478121029, 0.3691837478121029, 3, 1, 13.65974866666667, 40.979246,         0.383515069756994,  0.3835150697569941, 1.513573004259555, 13.622157038336, 0.1415875307953071, 0.1415875307953071),
    @("MuSCs", "Nppa", "Npr2", "FAPs",  2, 0.6666666666666666, 0.1108053333333333, 0.332416,           0.3691837478121029, 0.3691837478121029, 3, 1, 17.104695,         51.314085,         0.4802363832729211, 0.4802363832729211, 1.89529143104,     17.05762287936,  0.1772954678124265, 0.1772954678124265),
    @("MuSCs", "Nppa", "Npr2", "MuSCs", 2, 0.6666666666666666, 0.1108053333333333, 0.332416,           0.3691837478121029, 0.3691837478121029, 3, 1, 4.852797333333333, 14.558392,         0.1362485469700849, 0.1362485469700849, 0.5377158261191111, 4.839442435072, 0.05030074920436928, 0.05030074920436928)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
